$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (rows 2-29) holds a date serial that was bumped by one day
# (45394 -> 45395, i.e. 2024-04-12 -> 2024-04-13).
for ($row = 2; $row -le 29; $row++) {
    $ws.Cells.Item($row, 3).Value = 45395
}
